# edit.ps1
# Applies the "Updated cryptos list ... with GitHub Actions" data refresh:
#   - Updates Price (D) / Volume(1h) (E) columns for every existing row.
#   - Row 39/40 swap position: "FraxShare" <-> "Hedera" traded places in the
#     ranking, so columns B/C (name/link) are exchanged between the two rows
#     along with their own refreshed D/E values.
#
# Most of the new text values (percentages, multi-dot "thousands" prices,
# names, URLs) are unambiguous strings and Range.Value handles them as text
# automatically. A subset of the new Price values look exactly like plain
# decimal numbers (e.g. "1.007", "0.09827") - left alone, Excel's COM layer
# would auto-convert those to numeric cells. The source workbook stores
# every Price/Volume cell as text, so for that subset we briefly force a
# text NumberFormat while assigning, then restore the cell's original
# "Normal" style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Straightforward text assignments (names, links, percentages,
#      and prices that are unambiguously text already) ----
$ws.Range("D2").Value = "29.038.15"
$ws.Range("E2").Value = "  -4.25%  "
$ws.Range("D3").Value = "1.963.45"
$ws.Range("E3").Value = "  -6.23%  "
$ws.Range("E4").Value = "  +0.56%  "
$ws.Range("E5").Value = "  -4.26%  "
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("E7").Value = "  -6.12%  "
$ws.Range("E8").Value = "  -4.36%  "
$ws.Range("E9").Value = "  -2.42%  "
$ws.Range("E10").Value = "  -2.47%  "
$ws.Range("E11").Value = "  -6.84%  "
$ws.Range("E12").Value = "  -7.48%  "
$ws.Range("D13").Value = "2.008.21"
$ws.Range("E13").Value = "  +1.36%  "
$ws.Range("E14").Value = "  -8.33%  "
$ws.Range("E15").Value = "  -6.69%  "
$ws.Range("E16").Value = "  +0.57%  "
$ws.Range("E17").Value = "  -5.30%  "
$ws.Range("E18").Value = "  -10.32%  "
$ws.Range("E19").Value = "  -0.54%  "
$ws.Range("E20").Value = "  -9.26%  "
$ws.Range("E21").Value = "  +0.63%  "
$ws.Range("E22").Value = "  -6.11%  "
$ws.Range("D23").Value = "29.052.14"
$ws.Range("E23").Value = "  -4.14%  "
$ws.Range("E24").Value = "  -3.74%  "
$ws.Range("E25").Value = "  -1.05%  "
$ws.Range("D26").Value = "2.198.59"
$ws.Range("E26").Value = "  -4.79%  "
$ws.Range("E27").Value = "  -4.00%  "
$ws.Range("E28").Value = "  -5.85%  "
$ws.Range("E29").Value = "  -10.10%  "
$ws.Range("E30").Value = "  -9.51%  "
$ws.Range("E31").Value = "  -5.33%  "
$ws.Range("E32").Value = "  -8.76%  "
$ws.Range("E33").Value = "  -6.48%  "
$ws.Range("E34").Value = "  -8.80%  "
$ws.Range("E35").Value = "  -8.33%  "
$ws.Range("E36").Value = "  -5.89%  "
$ws.Range("E37").Value = "  -8.03%  "
$ws.Range("E38").Value = "  -3.23%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("E39").Value = "  -6.70%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("E40").Value = "  -11.97%  "
$ws.Range("E41").Value = "  -7.78%  "
$ws.Range("E42").Value = "  -9.74%  "
$ws.Range("E43").Value = "  -10.28%  "
$ws.Range("E44").Value = "  +0.43%  "
$ws.Range("E45").Value = "  -8.17%  "
$ws.Range("E46").Value = "  -6.65%  "
$ws.Range("E47").Value = "  -8.60%  "
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("E49").Value = "  -4.54%  "
$ws.Range("E50").Value = "  -2.92%  "
$ws.Range("E51").Value = "  -4.39%  "

# ---- Price values that look like plain numbers: must stay text ----
$forceTextCells = [ordered]@{
    "D5" = "327.38"
    "D7" = "0.4989"
    "D8" = "0.4203"
    "D9" = "52.72"
    "D10" = "0.09128"
    "D12" = "22.85"
    "D14" = "7.845"
    "D15" = "6.426"
    "D16" = "1.007"
    "D17" = "0.00001099"
    "D18" = "91.26"
    "D19" = "0.06672"
    "D22" = "5.958"
    "D24" = "12.04"
    "D27" = "156.20"
    "D28" = "20.56"
    "D29" = "6.169"
    "D30" = "2.258"
    "D31" = "126.69"
    "D33" = "0.09827"
    "D34" = "1.523"
    "D35" = "5.758"
    "D36" = "3.680"
    "D37" = "0.02413"
    "D38" = "1.299"
    "D39" = "0.06313"
    "D40" = "8.910"
    "D41" = "0.6430"
    "D42" = "11.41"
    "D43" = "0.1979"
    "D44" = "1.005"
    "D46" = "13.28"
    "D47" = "2.172"
    "D49" = "3.468"
    "D50" = "0.00000000332"
    "D51" = "0.06943"
}

foreach ($addr in $forceTextCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $forceTextCells[$addr]
    $cell.Style = "Normal"
}
